$d = $word.ActiveDocument

# 1. Replace "test" with "Préalable(s) : Aucun" in the Identification table
$d.Content.Find.Execute("test", $true, $false, $false, $false, $false, $true, 1, $false, "Préalable(s) : Aucun", 2)

$t2 = $d.Tables.Item(2)
$t2.Cell(2, 1).Range.Text = "Ceci est la présentation"

$t3 = $d.Tables.Item(3)
$t3.Cell(2, 1).Range.Text = "- Passer le cours facilement"

$t4 = $d.Tables.Item(4)
$t4.Cell(2, 1).Range.Text = "1er considération"

Write-Host ("Tables: " + $d.Tables.Count)
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    Write-Host ("Table $i rows: " + $t.Rows.Count + " cols: " + $t.Columns.Count)
    for ($r = 1; $r -le $t.Rows.Count; $r++) {
        for ($c = 1; $c -le $t.Columns.Count; $c++) {
            $cell = $t.Cell($r, $c)
            Write-Host ("  [$r,$c] = '" + $cell.Range.Text + "'")
        }
    }
}
